$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '66.272.71'
$ws.Range('E2').Value = '  -0.94%  '

# Row 3
$ws.Range('D3').Value = '3.534.42'
$ws.Range('E3').Value = '  +0.56%  '

# Row 4
$ws.Range('E4').Value = '  -0.09%  '

# Row 5
$ws.Range('D5').Value = '''607.32'
$ws.Range('E5').Value = '  -0.11%  '

# Row 6
$ws.Range('D6').Value = '''143.63'
$ws.Range('E6').Value = '  -3.01%  '

# Row 7
$ws.Range('D7').Value = '3.533.11'
$ws.Range('E7').Value = '  +0.57%  '

# Row 8
$ws.Range('E8').Value = '  -0.03%  '

# Row 11
$ws.Range('D11').Value = '''8.04'
$ws.Range('E11').Value = '  +1.05%  '

# Row 12
$ws.Range('E12').Value = '  -2.92%  '

# Row 13
$ws.Range('D13').Value = '4.129.74'
$ws.Range('E13').Value = '  +0.46%  '

# Row 14
$ws.Range('E14').Value = '  -4.52%  '

# Row 15
$ws.Range('D15').Value = '''30.19'
$ws.Range('E15').Value = '  -5.50%  '

# Row 16
$ws.Range('D16').Value = '3.529.35'
$ws.Range('E16').Value = '  +0.47%  '

# Row 17
$ws.Range('D17').Value = '66.339.34'
$ws.Range('E17').Value = '  -0.97%  '

# Row 18
$ws.Range('E18').Value = '  -0.70%  '

# Row 19
$ws.Range('D19').Value = '''10.96'
$ws.Range('E19').Value = '  +1.90%  '

# Row 20
$ws.Range('D20').Value = '''6.21'
$ws.Range('E20').Value = '  -3.91%  '

# Row 21
$ws.Range('D21').Value = '''14.92'
$ws.Range('E21').Value = '  -2.94%  '

# Row 22
$ws.Range('D22').Value = '''425.57'
$ws.Range('E22').Value = '  -2.93%  '

# Row 23
$ws.Range('D23').Value = '''0.602'
$ws.Range('E23').Value = '  -1.30%  '

# Row 24
$ws.Range('D24').Value = '''78.64'
$ws.Range('E24').Value = '  -1.04%  '

# Row 25
$ws.Range('D25').Value = '3.672.27'
$ws.Range('E25').Value = '  +0.53%  '

# Row 26
$ws.Range('E26').Value = '  +0.00%  '

# Row 27
$ws.Range('D27').Value = '''0.0000121'
$ws.Range('E27').Value = '  -1.01%  '

# Row 28
$ws.Range('D28').Value = '''8.08'
$ws.Range('E28').Value = '  -2.85%  '

# Row 29
$ws.Range('D29').Value = '''9.17'
$ws.Range('E29').Value = '  -6.21%  '

# Row 30
$ws.Range('E30').Value = '  -1.60%  '

# Row 31
$ws.Range('E31').Value = '  +0.26%  '

# Row 32
$ws.Range('B32').Value = 'Fetch.AI'
$ws.Range('C32').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D32').Value = '''1.49'
$ws.Range('E32').Value = '  -7.90%  '

# Row 33
$ws.Range('B33').Value = 'Kaspa'
$ws.Range('C33').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D33').Value = '''0.161'
$ws.Range('E33').Value = '  -4.42%  '

# Row 34
$ws.Range('E34').Value = '  -1.06%  '

# Row 35
$ws.Range('D35').Value = '3.521.59'
$ws.Range('E35').Value = '  +0.36%  '

# Row 36
$ws.Range('E36').Value = '  -0.02%  '

# Row 37
$ws.Range('D37').Value = '''1.75'
$ws.Range('E37').Value = '  -3.19%  '

# Row 38
$ws.Range('B38').Value = 'NEARProtocol'
$ws.Range('C38').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D38').Value = '''5.63'
$ws.Range('E38').Value = '  -5.68%  '

# Row 39
$ws.Range('B39').Value = 'Aptos'
$ws.Range('C39').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D39').Value = '''7.81'
$ws.Range('E39').Value = '  -2.96%  '

# Row 40
$ws.Range('E40').Value = '  -0.08%  '

# Row 41
$ws.Range('D41').Value = '''171.95'
$ws.Range('E41').Value = '  -1.24%  '

# Row 42
$ws.Range('D42').Value = '''0.0857'
$ws.Range('E42').Value = '  -4.33%  '

# Row 43
$ws.Range('D43').Value = '''5.19'
$ws.Range('E43').Value = '  -4.18%  '

# Row 44
$ws.Range('D44').Value = '''0.892'
$ws.Range('E44').Value = '  -0.31%  '

# Row 45
$ws.Range('E45').Value = '  -8.44%  '

# Row 46
$ws.Range('D46').Value = '''45.47'
$ws.Range('E46').Value = '  -1.52%  '

# Row 47
$ws.Range('E47').Value = '  -6.45%  '

# Row 48
$ws.Range('E48').Value = '  -5.09%  '

# Row 49
$ws.Range('D49').Value = '''2.41'
$ws.Range('E49').Value = '  -2.30%  '

# Row 50
$ws.Range('E50').Value = '  -4.38%  '

# Row 51
$ws.Range('D51').Value = '''0.946'
$ws.Range('E51').Value = '  -4.85%  '
